$p = $ppt.ActivePresentation

# --- Slide 2: merge the trailing " " + "3" runs of the "Lecture 3" title
#     paragraph into a single " 3" run (same visible text, fewer runs).
$titleShape = $p.Slides.Item(2).Shapes.Item(1)
$titleTextRange = $titleShape.TextFrame.TextRange
$lectureParagraph = $titleTextRange.Paragraphs(3, 1)
$tail = $lectureParagraph.Characters(8, 2)
$tail.Text = " 3"

# --- Remove the last slide (slide29.xml, sldId 285) from the deck.
$p.Slides.Item(29).Delete()
